# Auto-generated: apply cryptos.xlsx price/volume update (Sun Dec 31 2023 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text must remain stored as literal text (numeric-looking strings in the Price column),
# even though no value actually changes type: force text format, assign, then restore default style
# so the saved style index matches the untouched cells around them.
$textCells = @{
    'D2' = '42.588.56'
    'D3' = '2.294.70'
    'D5' = '316.15'
    'D6' = '103.97'
    'D7' = '0.624'
    'D9' = '0.601'
    'D10' = '39.38'
    'D11' = '0.0906'
    'D12' = '8.49'
    'D15' = '15.33'
    'D16' = '2.643.16'
    'D17' = '2.294.40'
    'D18' = '42.675.45'
    'D19' = '14.91'
    'D20' = '7.53'
    'D22' = '74.15'
    'D24' = '263.34'
    'D25' = '2.21'
    'D27' = '10.87'
    'D29' = '6.92'
    'D30' = '22.33'
    'D31' = '37.20'
    'D32' = '166.78'
    'D33' = '0.0874'
    'D34' = '0.131'
    'D35' = '2.61'
    'D37' = '4.58'
    'D38' = '0.0350'
    'D39' = '3.70'
    'D40' = '2.68'
    'D41' = '1.57'
    'D42' = '69.52'
    'D43' = '0.229'
    'D44' = '1.00'
    'D45' = '93.80'
    'D46' = '12.33'
    'D47' = '114.20'
    'D48' = '1.719.67'
    'D49' = '79.98'
    'D50' = '8.79'
    'D51' = '0.0997'
}
foreach ($ref in $textCells.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $textCells[$ref]
    $rng.Style = "Normal"
}

# Remaining cells (coin names, links, volume percentages) are safe to assign directly.
$plainCells = @{
    'E2' = '  -0.01%  '
    'E3' = '  -0.45%  '
    'E4' = '  -0.17%  '
    'E5' = '  -0.65%  '
    'E6' = '  +0.20%  '
    'E7' = '  -0.71%  '
    'E8' = '  -0.16%  '
    'E9' = '  -1.53%  '
    'E10' = '  -1.77%  '
    'E11' = '  -0.69%  '
    'E12' = '  +1.53%  '
    'E13' = '  +2.29%  '
    'E14' = '  +2.73%  '
    'E15' = '  -0.26%  '
    'E16' = '  -0.42%  '
    'E17' = '  -0.63%  '
    'E18' = '  +0.30%  '
    'E19' = '  +34.17%  '
    'E20' = '  -0.14%  '
    'E21' = '  -0.48%  '
    'E22' = '  +1.39%  '
    'E23' = '  -1.33%  '
    'E24' = '  -5.47%  '
    'E25' = '  -3.11%  '
    'E26' = '  +0.56%  '
    'E27' = '  +0.04%  '
    'E28' = '  -0.13%  '
    'E29' = '  +17.28%  '
    'E30' = '  -2.27%  '
    'E31' = '  +3.15%  '
    'E32' = '  +0.96%  '
    'E33' = '  -0.51%  '
    'E34' = '  -4.44%  '
    'E35' = '  +0.24%  '
    'E36' = '  -2.84%  '
    'E37' = '  -1.12%  '
    'E38' = '  -6.47%  '
    'E39' = '  -1.24%  '
    'E40' = '  -3.81%  '
    'E41' = '  +2.91%  '
    'E42' = '  -0.57%  '
    'E43' = '  +0.63%  '
    'B44' = 'FirstDigitalUSD'
    'C44' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'E44' = '  -0.16%  '
    'B45' = 'BitcoinSV'
    'C45' = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
    'E45' = '  -2.69%  '
    'E46' = '  +1.62%  '
    'E47' = '  +1.95%  '
    'B48' = 'Maker'
    'C48' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'E48' = '  +7.20%  '
    'B49' = 'ordi'
    'C49' = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
    'E49' = '  -0.73%  '
    'E50' = '  -1.80%  '
    'B51' = 'Cronos'
    'C51' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'E51' = '  -0.30%  '
}
foreach ($ref in $plainCells.Keys) {
    $ws.Range($ref).Value = $plainCells[$ref]
}
